$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.174468085106383
$ws.Range("C2").Value = 0.5957446808510638
$ws.Range("J2").Value = 0.00425531914893617
$ws.Range("P2").Value = 0.148936170212766
$ws.Range("S2").Value = 0.07659574468085106
$ws.Range("C3").Value = 0.03401360544217687
$ws.Range("J3").Value = 0.04081632653061224
$ws.Range("P3").Value = 0.7278911564625851
$ws.Range("S3").Value = 0.1972789115646258
$ws.Range("J4").Value = 0.1470588235294118
$ws.Range("O4").Value = 0.02941176470588235
$ws.Range("P4").Value = 0.6764705882352942
$ws.Range("S4").Value = 0.1470588235294118
$ws.Range("B6").Value = 0.05978260869565218
$ws.Range("D6").Value = 0.005434782608695652
$ws.Range("E6").Value = 0.005434782608695652
$ws.Range("F6").Value = 0.05978260869565218
$ws.Range("J6").Value = 0.2010869565217391
$ws.Range("O6").Value = 0.01630434782608696
$ws.Range("Q6").Value = 0.1956521739130435
$ws.Range("R6").Value = 0.05978260869565218
$ws.Range("S6").Value = 0.3967391304347826
$ws.Range("B7").Value = 0.1479289940828402
$ws.Range("D7").Value = 0.01775147928994083
$ws.Range("E7").Value = 0.005917159763313609
$ws.Range("F7").Value = 0.03550295857988166
$ws.Range("J7").Value = 0.1124260355029586
$ws.Range("O7").Value = 0.01775147928994083
$ws.Range("Q7").Value = 0.1656804733727811
$ws.Range("R7").Value = 0.04142011834319527
$ws.Range("S7").Value = 0.4556213017751479
$ws.Range("B8").Value = 0.07397959183673469
$ws.Range("D8").Value = 0.00510204081632653
$ws.Range("E8").Value = 0.002551020408163265
$ws.Range("F8").Value = 0.04336734693877551
$ws.Range("J8").Value = 0.1096938775510204
$ws.Range("O8").Value = 0.007653061224489796
$ws.Range("Q8").Value = 0.2066326530612245
$ws.Range("R8").Value = 0.09693877551020408
$ws.Range("S8").Value = 0.4540816326530612
$ws.Range("B9").Value = 0.07798165137614679
$ws.Range("D9").Value = 0.01834862385321101
$ws.Range("F9").Value = 0.07798165137614679
$ws.Range("J9").Value = 0.08256880733944955
$ws.Range("O9").Value = 0.009174311926605505
$ws.Range("Q9").Value = 0.2201834862385321
$ws.Range("R9").Value = 0.09174311926605505
$ws.Range("S9").Value = 0.4220183486238532
$ws.Range("B10").Value = 0.08783239323126511
$ws.Range("D10").Value = 0.0201450443190975
$ws.Range("F10").Value = 0.06365834004834811
$ws.Range("J10").Value = 0.1103948428686543
$ws.Range("O10").Value = 0.0145044319097502
$ws.Range("Q10").Value = 0.2159548751007252
$ws.Range("R10").Value = 0.07091055600322321
$ws.Range("S10").Value = 0.4165995165189363
$ws.Range("G11").Value = 0.1412213740458015
$ws.Range("J11").Value = 0.07633587786259542
$ws.Range("K11").Value = 0.2137404580152672
$ws.Range("L11").Value = 0.5496183206106871
$ws.Range("S11").Value = 0.01908396946564886
$ws.Range("F12").Value = 0.006666666666666667
$ws.Range("G12").Value = 0.74
$ws.Range("J12").Value = 0.1866666666666667
$ws.Range("L12").Value = 0.05333333333333334
$ws.Range("S12").Value = 0.01333333333333333
$ws.Range("G13").Value = 0.7714285714285715
$ws.Range("J13").Value = 0.2285714285714286
$ws.Range("G14").Value = 0.5
$ws.Range("J14").Value = 0.5
$ws.Range("F15").Value = 0.03191489361702127
$ws.Range("H15").Value = 0.1117021276595745
$ws.Range("I15").Value = 0.0851063829787234
$ws.Range("J15").Value = 0.3936170212765958
$ws.Range("K15").Value = 0.06382978723404255
$ws.Range("N15").Value = 0.005319148936170213
$ws.Range("O15").Value = 0.07446808510638298
$ws.Range("S15").Value = 0.2340425531914894
$ws.Range("F16").Value = 0.006211180124223602
$ws.Range("H16").Value = 0.124223602484472
$ws.Range("I16").Value = 0.1055900621118012
$ws.Range("J16").Value = 0.4782608695652174
$ws.Range("K16").Value = 0.09937888198757763
$ws.Range("M16").Value = 0.02484472049689441
$ws.Range("N16").Value = 0.006211180124223602
$ws.Range("O16").Value = 0.05590062111801242
$ws.Range("S16").Value = 0.09937888198757763
$ws.Range("F17").Value = 0.006451612903225806
$ws.Range("H17").Value = 0.1720430107526882
$ws.Range("I17").Value = 0.1225806451612903
$ws.Range("J17").Value = 0.4709677419354839
$ws.Range("K17").Value = 0.08602150537634409
$ws.Range("M17").Value = 0.01505376344086022
$ws.Range("O17").Value = 0.04086021505376344
$ws.Range("S17").Value = 0.08602150537634409
$ws.Range("F18").Value = 0.006172839506172839
$ws.Range("H18").Value = 0.1419753086419753
$ws.Range("I18").Value = 0.09259259259259259
$ws.Range("J18").Value = 0.4567901234567901
$ws.Range("K18").Value = 0.1049382716049383
$ws.Range("M18").Value = 0.01234567901234568
$ws.Range("O18").Value = 0.06172839506172839
$ws.Range("S18").Value = 0.1234567901234568
$ws.Range("F19").Value = 0.009876543209876543
$ws.Range("H19").Value = 0.1991769547325103
$ws.Range("I19").Value = 0.09382716049382717
$ws.Range("J19").Value = 0.405761316872428
$ws.Range("K19").Value = 0.09876543209876543
$ws.Range("M19").Value = 0.01810699588477366
$ws.Range("N19").Value = 0.001646090534979424
$ws.Range("O19").Value = 0.07325102880658436
$ws.Range("S19").Value = 0.09958847736625515
